# Update countries & provincias Spain
# Applies the data refresh captured in the target diff:
#  - swap the display order of Honduras/Barein (Honduras overtakes Barein in total cases)
#  - swap the display order of Groenlandia/Islas Malvinas
#  - bump the "Datos actualizados" timestamp
#  - refresh the numeric stats for Peru, Belgica, Honduras, Barein and Haiti rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 05:55"

# --- Row 10: Peru ---
$ws.Range("B10").Value = 384797
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 267850
$ws.Range("E10").Value = 98718
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 18229

# --- Row 37: Belgica ---
$ws.Range("B37").Value = 66026
$ws.Range("C37").Value = 299
$ws.Range("D37").Value = 17438
$ws.Range("E37").Value = 38767
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 9821

# --- Row 51: now Honduras (overtook Barein) ---
$ws.Range("A51").Value = "Honduras"
$ws.Range("B51").Value = 39276
$ws.Range("C51").Value = 838
$ws.Range("D51").Value = 4922
$ws.Range("E51").Value = 33238
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 18
$ws.Range("H51").Value = 1116

# --- Row 52: now Barein ---
$ws.Range("A52").Value = "Barein"
$ws.Range("B52").Value = 39131
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 35689
$ws.Range("E52").Value = 3302
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 140

# --- Row 91: Haiti ---
$ws.Range("B91").Value = 7315
$ws.Range("C91").Value = 18
$ws.Range("D91").Value = 4365
$ws.Range("E91").Value = 2793
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 157

# --- Row 210: now Groenlandia (was Islas Malvinas) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# --- Row 211: now Islas Malvinas (was Groenlandia) ---
$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
